# Auto-generated Excel COM-interop script applying the Spriggan_Profits market-data refresh.
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H:N) across all 8 job sheets
# with freshly pulled marketboard values from the scheduled runner.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC (82 cell updates) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1149.8
$ws.Range("I32").Value = 1083.3334
$ws.Range("J32").Value = 1249.5
$ws.Range("K32").Value = 1083.3334
$ws.Range("L32").Value = 1249.5
$ws.Range("M32").Value = -757.3334
$ws.Range("N32").Value = -1901.5
$ws.Range("H40").Value = 3484.375
$ws.Range("I40").Value = 2075
$ws.Range("J40").Value = 5833.3335
$ws.Range("K40").Value = 2075
$ws.Range("L40").Value = 5833.3335
$ws.Range("M40").Value = -1900
$ws.Range("N40").Value = -6183.3335
$ws.Range("H70").Value = 4074.5
$ws.Range("I70").Value = 1480.6923
$ws.Range("J70").Value = 7821.1113
$ws.Range("K70").Value = 4442.0769
$ws.Range("L70").Value = 23463.3339
$ws.Range("M70").Value = -4172.0769
$ws.Range("N70").Value = -24003.3339
$ws.Range("H73").Value = 4074.5
$ws.Range("I73").Value = 1480.6923
$ws.Range("J73").Value = 7821.1113
$ws.Range("K73").Value = 4442.0769
$ws.Range("L73").Value = 23463.3339
$ws.Range("M73").Value = -3506.0769
$ws.Range("N73").Value = -25335.3339
$ws.Range("H86").Value = 6473.6
$ws.Range("I86").Value = 7555.8887
$ws.Range("J86").Value = 4850.1665
$ws.Range("K86").Value = 7555.8887
$ws.Range("L86").Value = 4850.1665
$ws.Range("M86").Value = -6432.8887
$ws.Range("N86").Value = -7096.1665
$ws.Range("H88").Value = 16781920
$ws.Range("I88").Value = 37040144
$ws.Range("J88").Value = 2756994.2
$ws.Range("K88").Value = 37040144
$ws.Range("L88").Value = 2756994.2
$ws.Range("M88").Value = -37039738
$ws.Range("N88").Value = -2757806.2
$ws.Range("H89").Value = 6473.6
$ws.Range("I89").Value = 7555.8887
$ws.Range("J89").Value = 4850.1665
$ws.Range("K89").Value = 37779.4435
$ws.Range("L89").Value = 24250.8325
$ws.Range("M89").Value = -32163.4435
$ws.Range("N89").Value = -35482.8325
$ws.Range("H91").Value = 16781920
$ws.Range("I91").Value = 37040144
$ws.Range("J91").Value = 2756994.2
$ws.Range("K91").Value = 37040144
$ws.Range("L91").Value = 2756994.2
$ws.Range("M91").Value = -37038740
$ws.Range("N91").Value = -2759802.2
$ws.Range("H92").Value = 810.0476
$ws.Range("I92").Value = 761.1667
$ws.Range("J92").Value = 1103.3334
$ws.Range("K92").Value = 761.1667
$ws.Range("L92").Value = 1103.3334
$ws.Range("M92").Value = 486.8333
$ws.Range("N92").Value = -3599.3334
$ws.Range("H106").Value = 2420
$ws.Range("I106").Value = 2518.2727
$ws.Range("K106").Value = 2518.2727
$ws.Range("M106").Value = -1887.2727
$ws.Range("H107").Value = 1340.8334
$ws.Range("I107").Value = 1340.8334
$ws.Range("K107").Value = 1340.8334
$ws.Range("M107").Value = 579.1666
$ws.Range("H110").Value = 59999.668
$ws.Range("J110").Value = 59999.668
$ws.Range("L110").Value = 59999.668
$ws.Range("N110").Value = -68179.66800000001
$ws.Range("H137").Value = 2490.9285
$ws.Range("I137").Value = 1637.4
$ws.Range("J137").Value = 4624.75
$ws.Range("K137").Value = 4912.200000000001
$ws.Range("L137").Value = 13874.25
$ws.Range("M137").Value = -2362.200000000001
$ws.Range("N137").Value = -18974.25

# --- Sheet: ARM (34 cell updates) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5730.129
$ws.Range("I32").Value = 4421.8623
$ws.Range("K32").Value = 4421.8623
$ws.Range("M32").Value = -4134.8623
$ws.Range("H61").Value = 52634416
$ws.Range("I61").Value = 58825292
$ws.Range("J61").Value = 11999.5
$ws.Range("K61").Value = 58825292
$ws.Range("L61").Value = 11999.5
$ws.Range("M61").Value = -58825080
$ws.Range("N61").Value = -12423.5
$ws.Range("H80").Value = 50000
$ws.Range("J80").Value = 50000
$ws.Range("L80").Value = 50000
$ws.Range("N80").Value = -51996
$ws.Range("H83").Value = 50000
$ws.Range("J83").Value = 50000
$ws.Range("L83").Value = 150000
$ws.Range("N83").Value = -159984
$ws.Range("H112").Value = 73248
$ws.Range("J112").Value = 73248
$ws.Range("L112").Value = 73248
$ws.Range("N112").Value = -76202
$ws.Range("H132").Value = 8336036
$ws.Range("I132").Value = 9093039
$ws.Range("K132").Value = 27279117
$ws.Range("M132").Value = -27276587
$ws.Range("H136").Value = 52634416
$ws.Range("I136").Value = 58825292
$ws.Range("J136").Value = 11999.5
$ws.Range("K136").Value = 176475876
$ws.Range("L136").Value = 35998.5
$ws.Range("M136").Value = -176473326
$ws.Range("N136").Value = -41098.5

# --- Sheet: BSM (26 cell updates) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1376.1333
$ws.Range("I20").Value = 1247.4445
$ws.Range("J20").Value = 1569.1666
$ws.Range("K20").Value = 1247.4445
$ws.Range("L20").Value = 1569.1666
$ws.Range("M20").Value = -1000.4445
$ws.Range("N20").Value = -2063.1666
$ws.Range("H22").Value = 1729.45
$ws.Range("I22").Value = 1706.2667
$ws.Range("J22").Value = 1799
$ws.Range("K22").Value = 1706.2667
$ws.Range("L22").Value = 1799
$ws.Range("M22").Value = -1533.2667
$ws.Range("N22").Value = -2145
$ws.Range("H86").Value = 4084.8
$ws.Range("I86").Value = 4084.8
$ws.Range("K86").Value = 4084.8
$ws.Range("M86").Value = -2961.8
$ws.Range("H89").Value = 4084.8
$ws.Range("I89").Value = 4084.8
$ws.Range("K89").Value = 20424
$ws.Range("M89").Value = -14808
$ws.Range("H105").Value = 3801
$ws.Range("I105").Value = 3569.6667
$ws.Range("K105").Value = 3569.6667
$ws.Range("M105").Value = -1822.6667

# --- Sheet: CRP (11 cell updates) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H54").Value = 36666.332
$ws.Range("J54").Value = 36666.332
$ws.Range("L54").Value = 36666.332
$ws.Range("N54").Value = -37982.332
$ws.Range("H105").Value = 1139.8889
$ws.Range("I105").Value = 1196.4286
$ws.Range("J105").Value = 942
$ws.Range("K105").Value = 1196.4286
$ws.Range("L105").Value = 942
$ws.Range("M105").Value = 550.5714
$ws.Range("N105").Value = -4436

# --- Sheet: CUL (34 cell updates) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 62.466667
$ws.Range("J2").Value = 65.333336
$ws.Range("L2").Value = 392.000016
$ws.Range("N2").Value = -618.000016
$ws.Range("H23").Value = 442.0909
$ws.Range("I23").Value = 431
$ws.Range("K23").Value = 1293
$ws.Range("M23").Value = -1058
$ws.Range("H55").Value = 2425
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()
$ws.Range("H70").Value = 7599.8335
$ws.Range("I70").Value = 6119.8
$ws.Range("J70").Value = 15000
$ws.Range("K70").Value = 18359.4
$ws.Range("L70").Value = 45000
$ws.Range("M70").Value = -18044.4
$ws.Range("N70").Value = -45630
$ws.Range("H73").Value = 7599.8335
$ws.Range("I73").Value = 6119.8
$ws.Range("J73").Value = 15000
$ws.Range("K73").Value = 18359.4
$ws.Range("L73").Value = 45000
$ws.Range("M73").Value = -17267.4
$ws.Range("N73").Value = -47184
$ws.Range("H131").Value = 3166
$ws.Range("J131").Value = 4900
$ws.Range("L131").Value = 14700
$ws.Range("N131").Value = -24780
$ws.Range("H140").Value = 2990.889
$ws.Range("I140").Value = 2990.889
$ws.Range("K140").Value = 8972.667000000001
$ws.Range("M140").Value = -3792.667000000001

# --- Sheet: GSM (36 cell updates) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2740
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 2740
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 2740
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -4736
$ws.Range("H83").Value = 2740
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 2740
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 13700
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -23684
$ws.Range("H102").Value = 1985.2142
$ws.Range("I102").Value = 1952
$ws.Range("J102").Value = 2184.5
$ws.Range("K102").Value = 1952
$ws.Range("L102").Value = 2184.5
$ws.Range("M102").Value = -330
$ws.Range("N102").Value = -5428.5
$ws.Range("H113").Value = 66480.19
$ws.Range("I113").Value = 70578.87
$ws.Range("J113").Value = 5000
$ws.Range("K113").Value = 70578.87
$ws.Range("L113").Value = 5000
$ws.Range("M113").Value = -68408.87
$ws.Range("N113").Value = -9340
$ws.Range("H122").Value = 1945.3077
$ws.Range("I122").Value = 1367.3182
$ws.Range("K122").Value = 4101.9546
$ws.Range("M122").Value = -1651.9546
$ws.Range("H132").Value = 6581118
$ws.Range("I132").Value = 8335276.5
$ws.Range("K132").Value = 25005829.5
$ws.Range("M132").Value = -25003299.5

# --- Sheet: LTW (20 cell updates) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2572.6365
$ws.Range("I22").Value = 2572.6365
$ws.Range("K22").Value = 2572.6365
$ws.Range("M22").Value = -2277.6365
$ws.Range("H27").Value = 2572.6365
$ws.Range("I27").Value = 2572.6365
$ws.Range("K27").Value = 2572.6365
$ws.Range("M27").Value = -2465.6365
$ws.Range("H40").Value = 2905.6667
$ws.Range("I40").Value = 2905.6667
$ws.Range("K40").Value = 2905.6667
$ws.Range("M40").Value = -2769.6667
$ws.Range("H55").Value = 460.63635
$ws.Range("I55").Value = 220.875
$ws.Range("K55").Value = 220.875
$ws.Range("M55").Value = -47.875
$ws.Range("H136").Value = 2499.5
$ws.Range("I136").Value = 2000
$ws.Range("K136").Value = 6000
$ws.Range("M136").Value = -3450

# --- Sheet: WVR (7 cell updates) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 21740886
$ws.Range("I136").Value = 21740886
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 65222658
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -65220108
$ws.Range("N136").ClearContents()

